# Applies the "gh-pages output generated at 456a3b4" update to
# 江西-漫展信息.xlsx — refreshed "want-to-go"/price/cover numbers on the
# existing rows of the 展览 and 全部类型 sheets, plus a brand-new
# "吉安·WF无线次元新星动漫博览会" entry inserted in date order (2024-08-10)
# ahead of the 高安 entry, which (together with everything after it)
# shifts down by one row.

$wb = $excel.ActiveWorkbook

function Update-Sheet {
    param(
        [string]$SheetName,
        [int]$InsertRow,          # row number the new entry gets inserted at
        [hashtable]$SimpleEdits,  # cellRef -> new value, applied before the insert
        [hashtable]$NewRow        # column letter (A-I) -> value for the inserted row
    )

    $ws = $wb.Worksheets.Item($SheetName)
    Write-Host "Updating sheet" $SheetName "insert at row" $InsertRow

    # 1) Plain numeric / text refreshes on rows that are not affected by the
    #    later row insertion.
    foreach ($ref in $SimpleEdits.Keys) {
        $ws.Range($ref).Value = $SimpleEdits[$ref]
    }

    # 2) Insert a new row, pushing the 高安 row (and everything below it)
    #    down by one.
    $ws.Rows.Item($InsertRow).Insert()

    # Pick up the bold/bordered/centered formatting used by every column-A
    # cell (style index the row above already carries) so the new index
    # cell matches its neighbours instead of Excel's auto-generated
    # no-border variant.
    $ws.Range("A" + ($InsertRow - 1)).Copy()
    $ws.Range("A" + $InsertRow).PasteSpecial(-4122)

    # 3) Fill in the new row's data. Column B/E hold date-shaped text
    #    ("2024-08-10" / "2024.08.10 ...") that Excel would otherwise
    #    silently reinterpret as a date serial, so force text formatting
    #    on the date column first.
    $ws.Range("B" + $InsertRow).NumberFormat = "@"

    foreach ($col in @("A","B","C","D","E","F","G","H","I")) {
        if ($NewRow.ContainsKey($col)) {
            $ws.Range($col + $InsertRow).Value = $NewRow[$col]
        }
    }

    # 4) The A column is a manually maintained sequential row index
    #    (A2=1, A3=2, ...). Renumber everything from the inserted row to
    #    the last row so it stays sequential after the insert.
    $dims = $ws.UsedRange.Rows.Count
    $lastRow = $dims
    for ($r = $InsertRow + 1; $r -le $lastRow; $r++) {
        $prev = $ws.Range("A" + ($r - 1)).Value()
        $ws.Range("A" + $r).Value = $prev + 1
    }
}

$newRowData = @{
    "B" = "2024-08-10"
    "C" = "吉安·WF无线次元新星动漫博览会"
    "D" = "吉安南大道133号 吉安市全民健身中心"
    "E" = "2024.08.10 09:00-08.10 17:00"
    "F" = 0
    "G" = 45
    "H" = "https://show.bilibili.com/platform/detail.html?id=88023"
    "I" = "//i0.hdslb.com/bfs/openplatform/202406/f95zVAmw1718246635629.jpeg"
}

# ---- 展览 sheet: insert happens at row 41 -------------------------------
$sheet1Edits = @{
    "F2"  = 344
    "F3"  = 280
    "F4"  = 1250
    "G4"  = 58.5
    "I4"  = "//i1.hdslb.com/bfs/openplatform/202406/OEU3ijdb1719299094349.jpeg"
    "F6"  = 28
    "F10" = 3469
    "F11" = 127
    "F16" = 596
    "F17" = 91
    "F18" = 741
    "F21" = 57
    "F24" = 2616
    "F25" = 5129
    "F27" = 75
    "F29" = 2925
    "F31" = 2244
    "F34" = 81
    "F35" = 113
    "F37" = 311
}
Update-Sheet "展览" 41 $sheet1Edits $newRowData

# ---- 全部类型 sheet: same data, one row further down (insert at row 42) --
$sheet4Edits = @{
    "F2"  = 344
    "F3"  = 280
    "F4"  = 1250
    "G4"  = 58.5
    "I4"  = "//i1.hdslb.com/bfs/openplatform/202406/OEU3ijdb1719299094349.jpeg"
    "F6"  = 28
    "F10" = 3469
    "F11" = 127
    "F17" = 596
    "F18" = 91
    "F19" = 741
    "F22" = 57
    "F25" = 2616
    "F26" = 5129
    "F28" = 75
    "F30" = 2928
    "F32" = 2244
    "F35" = 81
    "F36" = 113
    "F38" = 311
}
Update-Sheet "全部类型" 42 $sheet4Edits $newRowData

Write-Host "Update complete"
